# Generate Report for Handoff
# Adds a newly-discovered source file
# "e883c2ac-cd28-4d30-bbba-c84bafe8d33d.md" to the localization status
# report (Overview / zh-cn / de-de sheets), inserting it just above the
# existing ".localization-config" bookkeeping row on every sheet.

$wb = $excel.ActiveWorkbook

$newFile        = "e883c2ac-cd28-4d30-bbba-c84bafe8d33d.md"
$zhXlf          = "e883c2ac-cd28-4d30-bbba-c84bafe8d33d.99754d798ebb015a64234e0a25891e72ca1b275d.zh-cn.xlf"
$deXlf          = "e883c2ac-cd28-4d30-bbba-c84bafe8d33d.99754d798ebb015a64234e0a25891e72ca1b275d.de-de.xlf"
$zhHandoffTime  = "2016-03-09 09:47:12"
$deHandoffTime  = "2016-03-09 09:47:15"
$epoch          = "0001-01-01 00:00:00"

$newFileCommit  = "1445f4636a1e231bfed6b1493257a64266eea94f"
$zhXlfCommit    = "8df4a860f0d4ed35d285a3e62bcf01b1fbe2a660"
$deXlfCommit    = "5718863c301a285bd990f34b734dc1f4bad6e9fa"

$newFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$newFileCommit/e2e/$newFile"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/$newFileCommit/.localization-config"
$zhXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhXlfCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf"
$deXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deXlfCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf"

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name) / B (zh-cn) / C (de-de)
# Row 3 becomes the new file; the old row 3 (.localization-config) slides
# down to row 4.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $newFileUrl, "", "", $newFile)
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $configUrl, "", "", ".localization-config")
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

# ---------------------------------------------------------------------
# Sheet "zh-cn": detail columns for the zh-cn handoff/handback cycle.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newFileUrl, "", "", $newFile)
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Range("D3").Value = $zhHandoffTime
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $configUrl, "", "", ".localization-config")
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

# ---------------------------------------------------------------------
# Sheet "de-de": detail columns for the de-de handoff/handback cycle.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newFileUrl, "", "", $newFile)
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $deXlfUrl, "", "", $deXlf)
$wsDe.Range("D3").Value = $deHandoffTime
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $configUrl, "", "", ".localization-config")
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"
